$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the cells that previously held string data (shared strings trimmed
# back to just the one used by E6).
$ws.Range("F7").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("G18").ClearContents()
$ws.Range("B19").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("E24").ClearContents()

# New data cells.
$ws.Range("F6").Value = 0
$ws.Range("F16").Value = 851

# Update the selection to match the post-edit state.
[void]$ws.Range("H12").Select()
